# Insert a new weekly record at row 74 (Macroferia Regional de Talca - Cilantro),
# pushing all existing data rows (old 74..101) down by one (new 75..102).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing rows 74..101 down to make room for the new record.
$ws.Rows.Item(74).Insert()

# Populate the newly inserted row 74 with the new weekly observation.
$ws.Range("A74").Value2 = 5
$ws.Range("B74").Value2 = "Macroferia Regional de Talca"
$ws.Range("C74").Value2 = "Maule"
$ws.Range("D74").Value2 = 45120
$ws.Range("E74").Value2 = 7
$ws.Range("F74").Value2 = 100112040
$ws.Range("G74").Value2 = "Cilantro"
$ws.Range("H74").Value2 = "Sin especificar"
$ws.Range("I74").Value2 = "Primera"
$ws.Range("J74").Value2 = 150
$ws.Range("K74").Value2 = 12000
$ws.Range("L74").Value2 = 12000
$ws.Range("M74").Value2 = 12000
$ws.Range("N74").Value2 = "`$/caja 36 atados"
$ws.Range("O74").Value2 = "Región Metropolitana"
$ws.Range("P74").Value2 = 333
$ws.Range("Q74").Value2 = 36
$ws.Range("R74").Value2 = "Hortaliza"

# Make sure the date cell keeps the date number format used by the rest of column D.
$ws.Range("D74").NumberFormat = $ws.Range("D75").NumberFormat
